# Updated GDP scaling factor for geoeng/DACD
# The model previously scaled global DAC potential down to a U.S.-only
# potential using a U.S. GDP share of world GDP. This updates the source
# data and labels to instead compute an EU GDP share of world GDP (using
# newer 2019 World Bank GDP figures), which ripples through the dependent
# formulas on the Data, DACD-potential, DACD-energyintensity and
# DACD-capex sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "About" sheet: insert a new "EU and World GDP" source block, right
#    after the existing "Electricity, Heat, and CapEx Data" source block
#    (which occupies rows 5-12). This pushes the old rows 14+ down by 5.
# ---------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Rows("14:18").Insert()

$wsAbout.Range("B14").Value = "EU and World GDP"
$wsAbout.Range("B15").Value = "World Bank"
$wsAbout.Range("B16").Value = "GDP (current US$) data"

# Match the formatting of the other "source block" headers/links already
# on this sheet (B5 = header style, B9/B11 = hyperlink style).
$wsAbout.Range("B5").Copy()
$wsAbout.Range("B14").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$wsAbout.Hyperlinks.Add($wsAbout.Range("B17"), "https://data.worldbank.org/indicator/NY.GDP.MKTP.CD")

$wsAbout.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2. "Data" sheet: replace the U.S. GDP / World GDP / US GDP share block
#    with EU GDP / World GDP / EU GDP share, using updated 2019 figures,
#    and rename the downstream "U.S. DAC potential" header to
#    "EU DAC potential". All of the formulas referencing these cells
#    (rows 78-79, 83-84 on this sheet, and the TREND() formulas on the
#    other sheets) recalculate automatically.
# ---------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("Data")

$wsData.Range("A72").Value = "EU GDP"
$wsData.Range("B72").Value = 15.625999999999999
$wsData.Range("C72").Value = "trillion USD"
$wsData.Range("D72").Value = 2019

$wsData.Range("A73").Value = "World GDP"
$wsData.Range("B73").Value = 87.799000000000007
$wsData.Range("C73").Value = "trillion USD"
$wsData.Range("D73").Value = 2019

$wsData.Range("A74").Value = "EU GDP share"

$wsData.Range("A76").Value = "EU DAC potential"

# ---------------------------------------------------------------------
# 3. Cosmetic: match the zoom level used across the sheets and leave the
#    "DACD-capex" sheet as the active one / active selection, as in the
#    saved workbook.
# ---------------------------------------------------------------------
$wsPotential = $wb.Worksheets.Item("DACD-potential")
$wsEnergyIntensity = $wb.Worksheets.Item("DACD-energyintensity")
$wsCapex = $wb.Worksheets.Item("DACD-capex")

$wsAbout.Activate()
$wsAbout.Range("B24").Select()
$excel.ActiveWindow.Zoom = 85

$wsData.Activate()
$wsData.Range("B22").Select()
$excel.ActiveWindow.Zoom = 85

$wsPotential.Activate()
$excel.ActiveWindow.Zoom = 85

$wsEnergyIntensity.Activate()
$wsEnergyIntensity.Range("B4").Select()
$excel.ActiveWindow.Zoom = 85

$wsCapex.Activate()
$wsCapex.Range("B2").Select()
$excel.ActiveWindow.Zoom = 85

Write-Host "Done: EU GDP share = $($wsData.Range('B74').Value2)"
